# 23/12, Macro Append Entry_A, FIFO formula in Entries
#
# The "Append Entry" macro pulled a fresh quote for ticker DSY.PA (row 59)
# and, FIFO-style, dropped the oldest Earnings Date stamp ("Oct 30, 2025")
# that had rolled out of the active window, while refreshing the latest
# prices / dividend yields returned by the data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# FIFO: the Earnings Date entry for DSY.PA (D59, "Oct 30, 2025") has aged out
# -- clear it so the stale shared string is dropped from the workbook.
$ws.Range("D59").ClearContents()

# Refresh of 1-Year Target Price / Dividend Yield (%) pulled by the data macro.
$ws.Range("B3").Value  = 20.5
$ws.Range("E3").Value  = 0.012

$ws.Range("E7").Value  = 0.0070999996

$ws.Range("B8").Value  = 456.79883
$ws.Range("E8").Value  = 0.0076

$ws.Range("E9").Value  = 0.0171

$ws.Range("E12").Value = 0.020299999

$ws.Range("E13").Value = 0.021

$ws.Range("E14").Value = 0.0339

$ws.Range("E16").Value = 0.013300001

$ws.Range("E20").Value = 0.0349

$ws.Range("E22").Value = 0.0072000003

$ws.Range("E23").Value = 0.0231

$ws.Range("E24").Value = 0.051799998

$ws.Range("B25").Value = 475.66858

$ws.Range("E26").Value = 0.0225

$ws.Range("B30").Value = 622.5141599999999

$ws.Range("E31").Value = 0.028299998

$ws.Range("E32").Value = 0.0095

$ws.Range("E35").Value = 0.0104

$ws.Range("E40").Value = 0.0249

$ws.Range("E41").Value = 0.0113

$ws.Range("E42").Value = 0.0229

$ws.Range("B44").Value = 399.151

$ws.Range("E45").Value = 0.0115

$ws.Range("B48").Value = 456.79883
$ws.Range("E48").Value = 0.0076

$ws.Range("E50").Value = 0.0467

$ws.Range("E51").Value = 0.04

$ws.Range("E53").Value = 0.0349

$ws.Range("E56").Value = 0.0648

$ws.Range("E57").Value = 0.0257

$ws.Range("B58").Value = 12.99842
$ws.Range("E58").Value = 0.0446

$ws.Range("E61").Value = 0.0178

$ws.Range("B63").Value = 37.93043
$ws.Range("E63").Value = 0.017

$ws.Range("E64").Value = 0.006

$ws.Range("B65").Value = 12.77778

$ws.Range("E66").Value = 0.0292

$ws.Range("E67").Value = 0.0539

$ws.Range("E68").Value = 0.0238

$ws.Range("E69").Value = 0.0682

$ws.Range("E71").Value = 0.010299999
